$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

$ws.Range("E27").Value = 5
$ws.Range("B28").Value = 12

$ws.Range("G30").Select()
